$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename occasion in row 7 (entry 6): "Ravens Game" -> "First Ravens Game"
$ws.Range("B7").Value = "First Ravens Game"

# Add picture popups (column H) for existing rows 2-30
$ws.Range("H2").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/10881686_10155030956530472_4441633013622423410_n.jpg?_nc_cat=103&oh=e58f2683328f99eb78da483e50a1f4ea&oe=5C161DA9"" style=""width:auto;height:50%;"">"
$ws.Range("H3").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/417358_171325126316477_1304069161_n.jpg?_nc_cat=109&oh=58e997770183716233a0ec215176cfef&oe=5C43AAA2"" style=""width:auto;height:75%;"">"
$ws.Range("H4").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/428226_10151383224825472_531323711_n.jpg?_nc_cat=101&oh=a7e47353b102b2243959b844b08cfba4&oe=5C180695"" style=""width:auto;height:75%;"">"
$ws.Range("H5").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/10885494_10155029400610472_1878221333012357816_n.jpg?_nc_cat=102&oh=85a5fa79dfce191fc20bdf623941de1b&oe=5C51D77A"" style=""width:auto;height:75%;"">"
$ws.Range("H6").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/14729164_10157628070480472_2549399823511176984_n.jpg?_nc_cat=100&oh=79324822dc2394d1094dae018064a42c&oe=5C5F0172"" style=""width:75%;height:auto;"">"
$ws.Range("H7").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/42432207_10161011613560472_2542234869564440576_o.jpg?_nc_cat=103&oh=6f1251f27b323f710f2e3de7e74c3a39&oe=5C62BE92"" style=""width:auto;height:60%;"">"
$ws.Range("H8").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/10805803_10154901235105472_7608660262787893426_n.jpg?_nc_cat=104&oh=78c247d666763496b1b9f18c3fcedd91&oe=5C60A1AF"" style=""width:75%;height:auto;"">"
$ws.Range("H9").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/10685582_10154732241745472_6024607142167467387_n.jpg?_nc_cat=102&oh=1494f94ecde57f9a7f6be35ae035d80e&oe=5C546BD2"" style=""width:75%;height:AUTO;"">"
$ws.Range("H10").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/19797_10155405886360472_1632482625524258935_n.jpg?_nc_cat=101&oh=75f0df8f337a9b1a048ef1678054ad8f&oe=5C16F224"" style=""width:75%;height:auto;"">"
$ws.Range("H11").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/535323_10156391993930472_1393915873010677955_n.jpg?_nc_cat=108&oh=b0a1847e9b80f847e9a2f22960a389af&oe=5C42B4E1"" style=""width:auto;height:75%;"">"
$ws.Range("H12").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/16728_10155171767845472_1736650083411801148_n.jpg?_nc_cat=107&oh=edf479e4053b597573fcf1765180020d&oe=5C189341"" style=""width:50%;height:auto;"">"
$ws.Range("H13").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/17799376_10158486995710472_1570977125754334957_n.jpg?_nc_cat=103&oh=3054d3dd9fe0be69ebc6a09cec5401eb&oe=5C484F56"" style=""width:75%;height:auto;"">"
$ws.Range("H14").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t31.0-8/12496376_10156554717550472_3394337341673659606_o.jpg?_nc_cat=106&oh=38f0da4f5a7127e4401480bad1f9d49f&oe=5C5D9CA8"" style=""width:auto;height:60%;"">"
$ws.Range("H15").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/17884023_10158492207155472_1286082213450447941_n.jpg?_nc_cat=105&oh=92ced4d4c9a36a466cf97e484202c356&oe=5C54F831"" style=""width:75%;height:auto;"">"
$ws.Range("H16").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/42723614_10161027530845472_2577854609603166208_o.jpg?_nc_cat=101&oh=a472bcf0eeae2686bfddcb43936b6b57&oe=5C51E650"" style=""width:75%;height:auto;"">"
$ws.Range("H17").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/528965_218234728292183_2005747330_n.jpg?_nc_cat=107&oh=1cd5a80551172cb30ae422bbe734a4f7&oe=5C5DDDE4"" style=""width:75%;height:auto;"">"
$ws.Range("H18").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/385763_231669023615420_1121748693_n.jpg?_nc_cat=103&oh=e1cb77eb8620f763ed3fc89f424a797d&oe=5C527BA1"" style=""width:75%;height:auto;"">"
$ws.Range("H19").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/1480765_10153573612775472_1070556733_n.jpg?_nc_cat=102&oh=cadc3be04d7742460a7dc484a890f883&oe=5C50E54E"" style=""width:auto;height:75%;"">"
$ws.Range("H20").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/1150418_10153118999170472_1431542147_n.jpg?_nc_cat=109&oh=d2b541447213917a362b8eb66a1d4aac&oe=5C458CEF"" style=""width:75%;height:auto;"">"
$ws.Range("H21").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t31.0-8/1074534_423811754401145_1162936417_o.jpg?_nc_cat=104&oh=5b56687fd9ca2cc6ec0aa47bf5ae7617&oe=5C4E08B0"" style=""width:75%;height:auto;"">"
$ws.Range("H22").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/62619_10152767099955472_1988822311_n.jpg?_nc_cat=101&oh=ae74fa09b643a3852b35813df1e394a9&oe=5C5F6944"" style=""width:auto;height:75%;"">"
$ws.Range("H23").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/574472_10152576450870472_454279356_n.jpg?_nc_cat=107&oh=e726de4458c62b90213d36c0ccfa0e1e&oe=5C62DAED"" style=""width:auto;height:75%;"">"
$ws.Range("H24").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t31.0-8/477295_385610101554644_562621180_o.jpg?_nc_cat=100&oh=e49fcf645102e9c6cca7c0550580db7d&oe=5C54B4B2"" style=""width:auto;height:70%;"">"
$ws.Range("H25").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/602766_339962592786062_1066806054_n.jpg?_nc_cat=103&oh=0e3c29c01785438bd16eff873ef152de&oe=5C5DA752"" style=""width:75%;height:auto;"">"
$ws.Range("H26").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t31.0-8/1979285_10153886727100472_1085309532_o.jpg?_nc_cat=107&oh=31e9fc95096a387c4b2d907136826ba2&oe=5C52D13D"" style=""width:auto;height:75%;"">"
$ws.Range("H27").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/18101_10152289472765472_1687701262_n.jpg?_nc_cat=101&oh=1bbd103f56c93587b117a0305416acab&oe=5C145450"" style=""width:75%;height:auto;"">"
$ws.Range("H28").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/10409453_10154891079400472_5213033156147697270_n.jpg?_nc_cat=110&oh=cff33a22a879e4fcb2dc070e1158cfc4&oe=5C4C8B49"" style=""width:auto;height:75%;"">"
$ws.Range("H29").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/1376413_10154800938795472_2721843062770153549_n.jpg?_nc_cat=107&oh=6e4362e57bd46e9ce514be8b345159fb&oe=5C5D351D"" style=""width:auto;height:75%"">"
$ws.Range("H30").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/10614194_10154544650605472_4099051271521523930_n.jpg?_nc_cat=110&oh=54995187b4c15656b162edffed059dbb&oe=5C19A26C"" style=""width:75%;height:auto;"">"

# Add new rows 31-33 with full data
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "Andrew Met Bay's Dad"
$ws.Range("C31").Value = "Fashion Island"
$ws.Range("D31").Value = "Relationship"
$ws.Range("E31").Value = "Jan 18 2015"
$ws.Range("F31").Value = 13.825656
$ws.Range("G31").Value = 100.678804
$ws.Range("H31").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/10920900_10155126879570472_506276529898579401_n.jpg?_nc_cat=103&oh=325bff7775418f093e37cd022126c77d&oe=5C555DA0"" style=""width:75%;height:auto;"">"

$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "Bay First Day of Work"
$ws.Range("C32").Value = "Gaysorn"
$ws.Range("D32").Value = "Business"
$ws.Range("E32").Value = "Mar 24 2012"
$ws.Range("F32").Value = 13.744959
$ws.Range("G32").Value = 100.540659
$ws.Range("H32").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/561353_192107217571601_153601725_n.jpg?_nc_cat=100&oh=3b4512804ead9977b733f6fa11ef3b93&oe=5C6082E1"" style=""width:75%;height:auto;"">"

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "First USA Halloween Together"
$ws.Range("C33").Value = "Millersville, MD"
$ws.Range("D33").Value = "Event"
$ws.Range("E33").Value = "Oct 31 2016"
$ws.Range("F33").Value = 39.054984
$ws.Range("G33").Value = -76.63964
$ws.Range("H33").Value = "<img src=""https://scontent-iad3-1.xx.fbcdn.net/v/t1.0-9/14900384_10157672410480472_3353697473017192986_n.jpg?_nc_cat=107&oh=c08e298bcf3088af7364d936af9b0703&oe=5C5452F6"" style=""width:auto;height:75%;"">"

# Update view state to match: scrolled down with active cell F37 selected (off-sheet placeholder from authoring)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("F37").Select()
